# Fix GitHub PDF rendering problem: re-crop "Picture 3" on slide 1 (the
# screenshot of the commit list) and shrink its frame to match the new
# (tighter) crop so the picture keeps its on-slide position/width.
#
# Equivalent OOXML effect:
#   <p:blipFill> ... <a:srcRect/> <a:stretch><a:fillRect/></a:stretch>
# becomes (crop the bottom ~54.33% of the source image off):
#   <p:blipFill> ... <a:srcRect b="30058"/> <a:stretch><a:fillRect/></a:stretch>
# and the picture's shape height shrinks from 144.625pt to 101.1535pt
# (width/position are untouched).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Picture 3")

# Crop the bottom of the picture: PowerPoint stores this as a percentage
# of the picture's native (uncropped) height in the <a:srcRect b="..."/>
# attribute; CropBottom is expressed in points of that native height.
# Native height here is 180.75pt, and we want b="30058" (30.058%), i.e.
# 0.30058 * 180.75 = 54.329835 points of crop.
$shp.PictureFormat.CropBottom = 54.329835

# Shrink the shape to the new (post-crop) picture height, 1284649 EMU,
# while leaving its width/left/top exactly as they were.
$shp.Height = 101.15346456692913
